$d = $word.ActiveDocument
$t = $d.Tables(1)
$t.Cell(1,1).Range.Text = "461×8="
$t.Cell(1,2).Range.Text = "403×9="
$t.Cell(1,3).Range.Text = "669×9="
$t.Cell(1,4).Range.Text = "870×7="
$t.Cell(1,5).Range.Text = "606×5="
$t.Cell(5,1).Range.Text = "976×6="
$t.Cell(5,2).Range.Text = "463×4="
$t.Cell(5,3).Range.Text = "826×9="
$t.Cell(5,4).Range.Text = "730×7="
$t.Cell(5,5).Range.Text = "322×3="
$t.Cell(10,1).Range.Text = "961×9="
$t.Cell(10,2).Range.Text = "685×2="
$t.Cell(10,3).Range.Text = "759×7="
$t.Cell(10,4).Range.Text = "647×2="
$t.Cell(10,5).Range.Text = "311×6="
$t.Cell(15,1).Range.Text = "925×7="
$t.Cell(15,2).Range.Text = "545×4="
$t.Cell(15,3).Range.Text = "922×7="
$t.Cell(15,4).Range.Text = "617×9="
$t.Cell(15,5).Range.Text = "527×4="
$t.Cell(20,1).Range.Text = "212×7="
$t.Cell(20,2).Range.Text = "434×6="
$t.Cell(20,3).Range.Text = "488×7="
$t.Cell(20,4).Range.Text = "257×2="
$t.Cell(20,5).Range.Text = "946×7="
